$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values that look like plain numbers (e.g. "257.71") are written with a
# leading apostrophe so Excel stores them as text (matching the original
# inline-string cell type) instead of silently converting them to numbers;
# the style is then reset to "Normal" so no stray number-format is left behind.

$ws.Range('D2').Value = "44.205.85"
$ws.Range('E2').Value = "  +3.75%  "
$ws.Range('D3').Value = "2.252.40"
$ws.Range('E3').Value = "  +2.58%  "
$ws.Range('E4').Value = "  +0.01%  "
$ws.Range('D5').Value = "'257.71"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "  +2.68%  "
$ws.Range('D6').Value = "'80.56"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "  +8.30%  "
$ws.Range('D7').Value = "'0.629"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "  +2.76%  "
$ws.Range('E8').Value = "  +0.10%  "
$ws.Range('E9').Value = "  +3.49%  "
$ws.Range('D10').Value = "'43.48"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "  +8.03%  "
$ws.Range('D11').Value = "'0.0934"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "  +1.76%  "
$ws.Range('E12').Value = "  +4.54%  "
$ws.Range('D14').Value = "2.588.79"
$ws.Range('E14').Value = "  +2.53%  "
$ws.Range('D15').Value = "'14.82"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "  +3.64%  "
$ws.Range('D16').Value = "2.265.26"
$ws.Range('E16').Value = "  +3.28%  "
$ws.Range('E17').Value = "  +2.08%  "
$ws.Range('D18').Value = "44.106.28"
$ws.Range('E18').Value = "  +3.83%  "
$ws.Range('E19').Value = "  +2.54%  "
$ws.Range('D20').Value = "'71.70"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "  +0.77%  "
$ws.Range('E21').Value = "  +3.01%  "
$ws.Range('D22').Value = "'2.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "  +10.02%  "
$ws.Range('D23').Value = "'235.50"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "  +3.03%  "
$ws.Range('D24').Value = "'9.46"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "  +0.83%  "
$ws.Range('E25').Value = "  +0.11%  "
$ws.Range('D26').Value = "'10.88"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "  +1.84%  "
$ws.Range('D27').Value = "'41.03"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "  +9.79%  "
$ws.Range('E28').Value = "  -0.48%  "
$ws.Range('E29').Value = "  +1.79%  "
$ws.Range('E30').Value = "  -0.62%  "
$ws.Range('D31').Value = "'173.40"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "  +2.35%  "
$ws.Range('D32').Value = "'20.70"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "  +3.02%  "
$ws.Range('D33').Value = "'0.0879"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "  +9.62%  "
$ws.Range('E34').Value = "  +3.72%  "
$ws.Range('E35').Value = "  +7.92%  "
$ws.Range('D36').Value = "'0.122"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "  +1.45%  "
$ws.Range('E37').Value = "  +12.38%  "
$ws.Range('D38').Value = "'4.55"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "  +4.70%  "
$ws.Range('D39').Value = "'13.15"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "  +7.75%  "
$ws.Range('D40').Value = "'2.93"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "  +22.43%  "
$ws.Range('E41').Value = "  +3.94%  "
$ws.Range('D42').Value = "'63.12"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "  +7.20%  "
$ws.Range('D43').Value = "'5.51"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "  +4.84%  "
$ws.Range('D44').Value = "'0.205"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "  +3.09%  "
$ws.Range('D45').Value = "'104.33"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "  +1.29%  "
$ws.Range('E46').Value = "  +1.51%  "
$ws.Range('D47').Value = "'0.0999"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "  +2.20%  "
$ws.Range('E48').Value = "  -2.85%  "
$ws.Range('E49').Value = "  +2.54%  "
$ws.Range('E50').Value = "  +2.60%  "
